$d = $word.ActiveDocument

$d.Content.Find.Execute("118÷5=23, 3", $true, $false, $false, $false, $false, $true, 1, $false, "521÷3=173, 2", 2) | Out-Null
$d.Content.Find.Execute("769÷2=384, 1", $true, $false, $false, $false, $false, $true, 1, $false, "669÷6=111, 3", 2) | Out-Null
$d.Content.Find.Execute("941÷5=188, 1", $true, $false, $false, $false, $false, $true, 1, $false, "708÷7=101, 1", 2) | Out-Null
$d.Content.Find.Execute("999÷8=124, 7", $true, $false, $false, $false, $false, $true, 1, $false, "967÷4=241, 3", 2) | Out-Null
$d.Content.Find.Execute("251÷4=62, 3", $true, $false, $false, $false, $false, $true, 1, $false, "255÷8=31, 7", 2) | Out-Null
$d.Content.Find.Execute("393÷7=56, 1", $true, $false, $false, $false, $false, $true, 1, $false, "788÷9=87, 5", 2) | Out-Null
$d.Content.Find.Execute("977÷2=488, 1", $true, $false, $false, $false, $false, $true, 1, $false, "853÷3=284, 1", 2) | Out-Null
$d.Content.Find.Execute("962÷7=137, 3", $true, $false, $false, $false, $false, $true, 1, $false, "110÷2=55, 0", 2) | Out-Null
$d.Content.Find.Execute("203÷5=40, 3", $true, $false, $false, $false, $false, $true, 1, $false, "465÷9=51, 6", 2) | Out-Null
$d.Content.Find.Execute("994÷6=165, 4", $true, $false, $false, $false, $false, $true, 1, $false, "573÷7=81, 6", 2) | Out-Null
$d.Content.Find.Execute("218÷2=109, 0", $true, $false, $false, $false, $false, $true, 1, $false, "755÷8=94, 3", 2) | Out-Null
$d.Content.Find.Execute("890÷2=445, 0", $true, $false, $false, $false, $false, $true, 1, $false, "117÷3=39, 0", 2) | Out-Null
$d.Content.Find.Execute("584÷4=146, 0", $true, $false, $false, $false, $false, $true, 1, $false, "581÷2=290, 1", 2) | Out-Null
$d.Content.Find.Execute("275÷9=30, 5", $true, $false, $false, $false, $false, $true, 1, $false, "198÷5=39, 3", 2) | Out-Null
$d.Content.Find.Execute("728÷4=182, 0", $true, $false, $false, $false, $false, $true, 1, $false, "101÷6=16, 5", 2) | Out-Null
$d.Content.Find.Execute("280÷5=56, 0", $true, $false, $false, $false, $false, $true, 1, $false, "497÷6=82, 5", 2) | Out-Null
$d.Content.Find.Execute("940÷8=117, 4", $true, $false, $false, $false, $false, $true, 1, $false, "494÷8=61, 6", 2) | Out-Null
$d.Content.Find.Execute("722÷4=180, 2", $true, $false, $false, $false, $false, $true, 1, $false, "167÷4=41, 3", 2) | Out-Null
$d.Content.Find.Execute("820÷2=410, 0", $true, $false, $false, $false, $false, $true, 1, $false, "305÷7=43, 4", 2) | Out-Null
$d.Content.Find.Execute("262÷7=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "997÷6=166, 1", 2) | Out-Null
$d.Content.Find.Execute("626÷8=78, 2", $true, $false, $false, $false, $false, $true, 1, $false, "645÷8=80, 5", 2) | Out-Null
$d.Content.Find.Execute("741÷4=185, 1", $true, $false, $false, $false, $false, $true, 1, $false, "214÷5=42, 4", 2) | Out-Null
$d.Content.Find.Execute("681÷7=97, 2", $true, $false, $false, $false, $false, $true, 1, $false, "249÷8=31, 1", 2) | Out-Null
$d.Content.Find.Execute("754÷5=150, 4", $true, $false, $false, $false, $false, $true, 1, $false, "194÷2=97, 0", 2) | Out-Null
$d.Content.Find.Execute("439÷2=219, 1", $true, $false, $false, $false, $false, $true, 1, $false, "190÷6=31, 4", 2) | Out-Null
